# Auto-update price data: prepend a new day's row (2026-01-09) above the
# existing top data row, pushing the historical rows down by one (same
# as inserting a new row 2 in the sheet and shifting 2..50 -> 3..51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above the current row 2 (first data row, right
# below the header). All existing data rows shift down by one.
$ws.Rows.Item(2).Insert()

# The newly inserted row inherits formatting from the row above (the
# bold/bordered header); strip that so it matches the plain data rows.
$ws.Range("A2:D2").ClearFormats()

# Populate the new row with today's price data. Force column A to be
# treated as plain text (matching the other date cells, which are
# strings rather than real date values) before assigning the value.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-01-09"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# Drop the number-format override now that the value is safely stored
# as text, so the cell ends up with no explicit style (like the rest
# of the data rows).
$ws.Range("A2:D2").ClearFormats()
